# Update Nalco PDF (2025-08-15 04:03:20 UTC)
#
# The price sheet gets a brand-new top data row (Sl.no. 2, the 15-Aug-2025
# circular) and the previous top row (Sl.no. 1, the 07-Aug-2025 circular)
# shifts down to row 3, unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldUrl = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"
$newUrl = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"

# 1) Push the existing data row (row 2) down to row 3 by inserting a blank
#    row above it.
$ws.Rows(2).Insert()

# 2) The inserted row has no formatting of its own yet - clone it from the
#    row it pushed down (row 3), which still carries the original look
#    (centered text, "0.000" number format on the price column, etc.).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Fill the new row with the latest circular's figures.
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "ALUMINIUM INGOT"
$ws.Range("C2").Value = "IE07"
$ws.Range("D2").Value = 269.45
$ws.Range("E2").Value = "15-08-2025"
$ws.Range("F2").Value = $newUrl

# 4) Rebuild the hyperlinks: the new circular link lives on F2, and the
#    older circular (now on row 3) keeps its own link pointing at the
#    original PDF.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), $newUrl)
$ws.Hyperlinks.Add($ws.Range("F3"), $oldUrl)

# Adding a hyperlink re-styles its cell with the blue/underlined built-in
# "Hyperlink" look; put the plain centered style back so F2/F3 match the
# rest of their rows, exactly like the source sheet.
$ws.Range("C2").Copy()
$ws.Range("F2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
